$wb = $excel.ActiveWorkbook


# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 5047.485
$ws.Range("J17").Value2 = 5145.875
$ws.Range("L17").Value2 = 15437.625
$ws.Range("N17").Value2 = -15773.625
$ws.Range("H48").Value2 = 2906.3333
$ws.Range("J48").Value2 = 2906.3333
$ws.Range("L48").Value2 = 8718.999899999999
$ws.Range("N48").Value2 = -9302.999899999999
$ws.Range("H56").Value2 = 2906.3333
$ws.Range("J56").Value2 = 2906.3333
$ws.Range("L56").Value2 = 8718.999899999999
$ws.Range("N56").Value2 = -9786.999899999999
$ws.Range("H87").Value2 = 39527
$ws.Range("J87").Value2 = 39527
$ws.Range("L87").Value2 = 39527
$ws.Range("N87").Value2 = -42023
$ws.Range("H90").Value2 = 39527
$ws.Range("J90").Value2 = 39527
$ws.Range("L90").Value2 = 118581
$ws.Range("N90").Value2 = -131061
$ws.Range("H100").Value2 = 2547.7856
$ws.Range("I100").Value2 = 1709.8572
$ws.Range("J100").Value2 = 3385.7144
$ws.Range("K100").Value2 = 1709.8572
$ws.Range("L100").Value2 = 3385.7144
$ws.Range("M100").Value2 = -1168.8572
$ws.Range("N100").Value2 = -4467.7144
$ws.Range("H129").Value2 = 1520.8055
$ws.Range("J129").Value2 = 1640.9697
$ws.Range("L129").Value2 = 4922.909100000001
$ws.Range("N129").Value2 = -14922.9091
$ws.Range("H137").Value2 = 112987.695
$ws.Range("I137").Value2 = 139539.89
$ws.Range("J137").Value2 = 2985.7144
$ws.Range("K137").Value2 = 418619.67
$ws.Range("L137").Value2 = 8957.143199999999
$ws.Range("M137").Value2 = -416069.67
$ws.Range("N137").Value2 = -14057.1432
$ws.Range("H138").Value2 = 4029.551
$ws.Range("I138").Value2 = 3542.8572
$ws.Range("J138").Value2 = 4110.6665
$ws.Range("K138").Value2 = 10628.5716
$ws.Range("L138").Value2 = 12331.9995
$ws.Range("M138").Value2 = -5488.571599999999
$ws.Range("N138").Value2 = -22611.9995

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1070.9048
$ws.Range("I2").Value2 = 1024.3125
$ws.Range("J2").Value2 = 1220
$ws.Range("K2").Value2 = 1024.3125
$ws.Range("L2").Value2 = 1220
$ws.Range("M2").Value2 = -911.3125
$ws.Range("N2").Value2 = -1446
$ws.Range("H32").Value2 = 8872.781000000001
$ws.Range("I32").Value2 = 6502.349
$ws.Range("K32").Value2 = 6502.349
$ws.Range("M32").Value2 = -6215.349
$ws.Range("H116").Value2 = 1070.9048
$ws.Range("I116").Value2 = 1024.3125
$ws.Range("J116").Value2 = 1220
$ws.Range("K116").Value2 = 1024.3125
$ws.Range("L116").Value2 = 1220
$ws.Range("M116").Value2 = 1269.6875
$ws.Range("N116").Value2 = -5808
$ws.Range("H122").Value2 = 1778.2927
$ws.Range("I122").Value2 = 1682.9429
$ws.Range("K122").Value2 = 5048.8287
$ws.Range("M122").Value2 = -2598.8287

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1070.9048
$ws.Range("I3").Value2 = 1024.3125
$ws.Range("J3").Value2 = 1220
$ws.Range("K3").Value2 = 1024.3125
$ws.Range("L3").Value2 = 1220
$ws.Range("M3").Value2 = -910.3125
$ws.Range("N3").Value2 = -1448
$ws.Range("H35").Value2 = 24988
$ws.Range("J35").Value2 = 24988
$ws.Range("L35").Value2 = 24988
$ws.Range("N35").Value2 = -25608
$ws.Range("H134").Value2 = 3519.2559
$ws.Range("I134").Value2 = 3555.4285
$ws.Range("J134").Value2 = 2000
$ws.Range("K134").Value2 = 10666.2855
$ws.Range("L134").Value2 = 6000
$ws.Range("M134").Value2 = -8131.2855
$ws.Range("N134").Value2 = -11070

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3706.5574
$ws.Range("I31").Value2 = 1883.8387
$ws.Range("K31").Value2 = 1883.8387
$ws.Range("M31").Value2 = -1588.8387
$ws.Range("H34").Value2 = 3706.5574
$ws.Range("I34").Value2 = 1883.8387
$ws.Range("K34").Value2 = 1883.8387
$ws.Range("M34").Value2 = -1681.8387
$ws.Range("H58").Value2 = 16124.206
$ws.Range("I58").Value2 = 1372.4445
$ws.Range("J58").Value2 = 32719.938
$ws.Range("K58").Value2 = 1372.4445
$ws.Range("L58").Value2 = 32719.938
$ws.Range("M58").Value2 = -1169.4445
$ws.Range("N58").Value2 = -33125.93799999999
$ws.Range("H107").Value2 = 1274.2727
$ws.Range("I107").Value2 = 456.8889
$ws.Range("K107").Value2 = 456.8889
$ws.Range("M107").Value2 = 1463.1111
$ws.Range("H136").Value2 = 16124.206
$ws.Range("I136").Value2 = 1372.4445
$ws.Range("J136").Value2 = 32719.938
$ws.Range("K136").Value2 = 4117.333500000001
$ws.Range("L136").Value2 = 98159.814
$ws.Range("M136").Value2 = -1567.333500000001
$ws.Range("N136").Value2 = -103259.814
$ws.Range("H137").Value2 = 25000
$ws.Range("J137").Value2 = 25000
$ws.Range("L137").Value2 = 25000
$ws.Range("N137").Value2 = -35200

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value2 = 0
$ws.Range("J82").Value2 = 0
$ws.Range("L82").Value2 = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value2 = 0
$ws.Range("J85").Value2 = 0
$ws.Range("L85").Value2 = 0
$ws.Range("N85").ClearContents()
$ws.Range("H107").Value2 = 6869.6206
$ws.Range("I107").Value2 = 8522.956
$ws.Range("J107").Value2 = 531.8333
$ws.Range("K107").Value2 = 25568.868
$ws.Range("L107").Value2 = 1595.4999
$ws.Range("M107").Value2 = -23648.868
$ws.Range("N107").Value2 = -5435.4999
$ws.Range("H113").Value2 = 992.8570999999999
$ws.Range("I113").Value2 = 0
$ws.Range("J113").Value2 = 992.8570999999999
$ws.Range("K113").Value2 = 0
$ws.Range("L113").Value2 = 2978.5713
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value2 = -7318.5713
$ws.Range("H129").Value2 = 224130.83
$ws.Range("J129").Value2 = 287978.22
$ws.Range("L129").Value2 = 863934.6599999999
$ws.Range("N129").Value2 = -873934.6599999999
$ws.Range("H130").Value2 = 2100
$ws.Range("J130").Value2 = 2200
$ws.Range("L130").Value2 = 6600
$ws.Range("N130").Value2 = -16640
$ws.Range("H131").Value2 = 721.23
$ws.Range("J131").Value2 = 754.9231
$ws.Range("L131").Value2 = 2264.7693
$ws.Range("N131").Value2 = -12344.7693
$ws.Range("H137").Value2 = 14497116
$ws.Range("I137").Value2 = 1317.2727
$ws.Range("J137").Value2 = 27784930
$ws.Range("K137").Value2 = 3951.8181
$ws.Range("L137").Value2 = 83354790
$ws.Range("M137").Value2 = 1148.1819
$ws.Range("N137").Value2 = -83364990

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 26319384
$ws.Range("I80").Value2 = 41669864
$ws.Range("K80").Value2 = 41669864
$ws.Range("M80").Value2 = -41668866
$ws.Range("H83").Value2 = 26319384
$ws.Range("I83").Value2 = 41669864
$ws.Range("K83").Value2 = 208349320
$ws.Range("M83").Value2 = -208344328
$ws.Range("H107").Value2 = 356.73334
$ws.Range("I107").Value2 = 331.22223
$ws.Range("J107").Value2 = 395
$ws.Range("K107").Value2 = 331.22223
$ws.Range("L107").Value2 = 395
$ws.Range("M107").Value2 = 1588.77777
$ws.Range("N107").Value2 = -4235
$ws.Range("H113").Value2 = 6011.8423
$ws.Range("I113").Value2 = 8894.091
$ws.Range("J113").Value2 = 2048.75
$ws.Range("K113").Value2 = 8894.091
$ws.Range("L113").Value2 = 2048.75
$ws.Range("M113").Value2 = -6724.091
$ws.Range("N113").Value2 = -6388.75
$ws.Range("H132").Value2 = 104960.12
$ws.Range("I132").Value2 = 139046.06
$ws.Range("J132").Value2 = 53831.2
$ws.Range("K132").Value2 = 417138.18
$ws.Range("L132").Value2 = 161493.6
$ws.Range("M132").Value2 = -414608.18
$ws.Range("N132").Value2 = -166553.6

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 2050
$ws.Range("I136").Value2 = 1957.8125
$ws.Range("J136").Value2 = 5000
$ws.Range("K136").Value2 = 5873.4375
$ws.Range("L136").Value2 = 15000
$ws.Range("M136").Value2 = -3323.4375
$ws.Range("N136").Value2 = -20100

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 428.57144
$ws.Range("I100").Value2 = 400
$ws.Range("J100").Value2 = 600
$ws.Range("K100").Value2 = 800
$ws.Range("L100").Value2 = 1200
$ws.Range("M100").Value2 = -259
$ws.Range("N100").Value2 = -2282
